$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Second paragraph currently reads:
#   "But, what are its semantics, and are those semantics useful?"
$para = $tr.Paragraphs(2)

$part1 = "But its behavior is not really " + [char]0x201C + "all-or-none" + [char]0x201D + ".  It is all, or none, or any prefix.  This is the "
$part2 = "same"
$part3 = " set of final resulting states as you can be in when you "
$part4 = "use the "
$part5 = [char]0x201C + "stop on first error" + [char]0x201D + " mode"

$fullText = $part1 + $part2 + $part3 + $part4 + $part5

# Replace the paragraph's text in one shot. PowerPoint's autocorrect-style
# diffing preserves a literal "But" prefix match against the old text as
# its own run; normalize back to a single run by re-assigning the full
# span to itself before carving out the runs we actually want.
$para.Text = $fullText
$whole = $para.Characters(1, $fullText.Length)
$whole.Text = $fullText

$start = 1
$r1 = $para.Characters($start, $part1.Length)
$start += $part1.Length

$r2 = $para.Characters($start, $part2.Length)
$r2.Font.Underline = -1
$start += $part2.Length

$r3 = $para.Characters($start, $part3.Length)
$start += $part3.Length

$r4 = $para.Characters($start, $part4.Length)
$start += $part4.Length

$r5 = $para.Characters($start, $part5.Length)
$start += $part5.Length

# r3/r4/r5 share identical (no) direct formatting, so they'd otherwise
# collapse back into one run; re-assert each span's own text (from the
# original literals, not read back, since Text reads normalize curly
# quotes to straight ones) to force PowerPoint to keep them as separate
# runs.
$r3.Text = $part3
$r4.Text = $part4
$r5.Text = $part5
